$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared / rich text cells (header strings) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Structural change: text -> number (set NumberFormat then Value) ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -50
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 3
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1

# --- Structural change: number -> text (copy format+value from donor "0" text cell) ---
$ws.Range("C28").Copy($ws.Range("F28"))
$ws.Range("C29").Copy($ws.Range("F29"))
$ws.Range("C30").Copy($ws.Range("F30"))

# --- Plain numeric value updates ---
# Row 14
$ws.Range("N14").Value = -92
# Row 15
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = -28.571428571428
$ws.Range("L15").Value = -25
$ws.Range("M15").Value = -16.666666666666
$ws.Range("N15").Value = -65.116279069767
# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 80
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 18.75
$ws.Range("I16").Value = 209
$ws.Range("J16").Value = 192
$ws.Range("K16").Value = 8.854166666666
$ws.Range("L16").Value = 50.359712230215
$ws.Range("M16").Value = -16.733067729083
$ws.Range("N16").Value = -76.543209876543
# Row 17
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -15.384615384615
$ws.Range("I17").Value = 278
$ws.Range("J17").Value = 223
$ws.Range("K17").Value = 24.663677130044
$ws.Range("L17").Value = 34.951456310679
$ws.Range("M17").Value = 24.107142857142
$ws.Range("N17").Value = -47.842401500938
# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("F18").Value = 25
$ws.Range("H18").Value = -3.846153846153
$ws.Range("I18").Value = 142
$ws.Range("J18").Value = 186
$ws.Range("K18").Value = -23.655913978494
$ws.Range("L18").Value = -9.554140127388
$ws.Range("M18").Value = -37.991266375545
$ws.Range("N18").Value = -82.316313823163
# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -43.75
$ws.Range("F19").Value = 54
$ws.Range("H19").Value = -20.588235294117
$ws.Range("I19").Value = 438
$ws.Range("J19").Value = 427
$ws.Range("K19").Value = 2.576112412177
$ws.Range("L19").Value = 45.033112582781
$ws.Range("M19").Value = 140.659340659341
$ws.Range("N19").Value = 27.696793002915
# Row 20
$ws.Range("C20").Value = 7
$ws.Range("E20").Value = -30
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 127
$ws.Range("J20").Value = 152
$ws.Range("K20").Value = -16.447368421052
$ws.Range("L20").Value = 29.591836734693
$ws.Range("M20").Value = 33.684210526315
$ws.Range("N20").Value = -78.868552412645
# Row 21
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = -17.647058823529
$ws.Range("F21").Value = 174
$ws.Range("G21").Value = 201
$ws.Range("H21").Value = -13.432835820895
$ws.Range("I21").Value = 1211
$ws.Range("J21").Value = 1203
$ws.Range("K21").Value = 0.665004156275
$ws.Range("L21").Value = 30.636461704422
$ws.Range("M21").Value = 20.138888888888
$ws.Range("N21").Value = -62.611917258413
# Row 22
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = 137.5
$ws.Range("M22").Value = 35.714285714285
# Row 23
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 8
$ws.Range("J23").Value = 21
$ws.Range("K23").Value = 4.761904761904
$ws.Range("M23").Value = 214.285714285714
# Row 24
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -4.347826086956
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 628
$ws.Range("J24").Value = 615
$ws.Range("K24").Value = 2.113821138211
$ws.Range("L24").Value = 12.142857142857
$ws.Range("M24").Value = 41.760722347629
# Row 25
$ws.Range("C25").Value = 11
$ws.Range("E25").Value = -15.384615384615
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = 10.63829787234
$ws.Range("I25").Value = 415
$ws.Range("J25").Value = 380
$ws.Range("K25").Value = 9.210526315789
$ws.Range("L25").Value = 32.165605095541
$ws.Range("M25").Value = -16.498993963782
# Row 26
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 28
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = -6.666666666666
$ws.Range("L26").Value = -9.677419354838
# Row 27
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 35
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = 16.666666666666
# Row 28
$ws.Range("H28").Value = -100
$ws.Range("L28").Value = -76.666666666666
$ws.Range("N28").Value = -94.85294117647
# Row 29
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -70
$ws.Range("N29").Value = -95.348837209302
